$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alphabet soup")

# Row 2 (ce_as_002): fill in the EXPECTED RESULT (column D) with the search results list
$ws.Range("D2").Value = "['Palabra AIRE encontrada entre (0, 0) y (3, 3).', 'Palabra AGUA encontrada entre (0, 0) y (0, 3).', 'Palabra TIERRA no encontrada.', 'Palabra NIU encontrada entre (2, 0) y (0, 2).']"

# Row 3 (ce_as_003): DESCRIPTION (column C) changes from "no separator" to "contains null word (length 1)"
$ws.Range("C3").Value = "contains null word (length 1)"

# Row 5 (ce_as_005): EXPECTED RESULT (column D) filled in
$ws.Range("D5").Value = "Incorrect table format."

# Row 7 (ce_as_006): VALIDITY/DESCRIPTION/EXPECTED RESULT filled in
$ws.Range("B7").Value = "NO"
$ws.Range("C7").Value = "non alphabetical characters in table"
$ws.Range("D7").Value = "Incorrect format. Game file corrupted."

# Row 8 (ce_as_007): VALIDITY/DESCRIPTION/EXPECTED RESULT filled in
$ws.Range("B8").Value = "NO"
$ws.Range("C8").Value = "non alphabetical characters in wordlist"
$ws.Range("D8").Value = "Incorrect format. Game file corrupted."

$wb.Save()
